$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 2998
$ws.Range("C3").Value = 2881
$ws.Range("C4").Value = 2123
$ws.Range("C5").Value = 1304
$ws.Range("C6").Value = 1178
$ws.Range("C7").Value = 701
$ws.Range("C8").Value = 598
$ws.Range("C9").Value = 438
$ws.Range("C10").Value = 434

$ws.Range("A11").Value = "Textiles & Cozy Items"
$ws.Range("B11").Value = "Textiles & Cozy Items"
$ws.Range("C11").Value = 409
